$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.629.39'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.598.37'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.515'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.80%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0619'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("E9").Value = '  +0.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.51'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.71%  '
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.822.73'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.594.14'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("E15").Value = '  -0.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.79'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.616.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("E18").Value = '  +1.11%  '
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.27'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.00%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.49'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.77%  '
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.50%  '
$ws.Range("E28").Value = '  -0.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.31'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0511'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.54%  '
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.23'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.95'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.277.24'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.618'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.45'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("E37").Value = '  +0.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0171'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.28%  '
$ws.Range("B39").Value = 'WEMIXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.07'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +19.31%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.838'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.48'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.58%  '
$ws.Range("E42").Value = '  +0.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.785'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.06'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.736.01'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.05'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.60'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.82%  '
$ws.Range("E48").Value = '  +3.92%  '
$ws.Range("E49").Value = '  +0.87%  '
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("E51").Value = '  -0.96%  '
